# Daily update at 8 AM UTC
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The row that was previously the most recent day (row 22) is no longer the
# latest, so it reverts from the "latest day" date-only format back to the
# standard datetime format shared by all the other historical rows.
$ws.Range("A22").NumberFormat = $ws.Range("A21").NumberFormat

# Append today's data as the new latest row (row 23), using the date-only
# format that marks the most-recent day.
$ws.Range("A23").Value = 45972
$ws.Range("A23").NumberFormat = "YYYY-MM-DD"
$ws.Range("B23").Value = 49
$ws.Range("C23").Value = 57
$ws.Range("D23").Value = 56
